# Daily attendance processing - 2025-11-03 08:54:48
# Reproduces the recorded-by reordering, status/count corrections and the
# "Pending" -> "Not Recorded" recoloring for the three sessions that rolled
# past their date without being recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Recorded By" text reordering: "dnasr281@gmail.com, System"
#    -> "System, dnasr281@gmail.com" (81 sessions across all groups)
# ---------------------------------------------------------------------
$swapRows = @(3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,90,92,93,94,96,99,101,109,110,111,112,116,118,119,120,122,125,127,135,136,137,138,142,144,145,146,148,151,153)
foreach ($r in $swapRows) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}

# Similar reordering for the "backup@backdoor.com, System, system" rows
$backupRows = @(2,28,54)
foreach ($r in $backupRows) {
    $ws.Cells.Item($r, 7).Value = "backup@backdoor.com, system, System"
}

# And the "dnasr281@gmail.com, admin@admin.com" rows
$adminRows = @(87,113,139)
foreach ($r in $adminRows) {
    $ws.Cells.Item($r, 7).Value = "admin@admin.com, dnasr281@gmail.com"
}

# ---------------------------------------------------------------------
# 2) Class statistics correction (B2A block): Missing/Pending sessions
# ---------------------------------------------------------------------
$ws.Range("L7").Value = 3   # Missing Sessions
$ws.Range("L8").Value = 9   # Pending Sessions

# ---------------------------------------------------------------------
# 3) Group statistics: Missing / Pending counts swapped for B2A, B2B, B2C
# ---------------------------------------------------------------------
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 0
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 0
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 0

# ---------------------------------------------------------------------
# 4) Sessions that are now overdue ("Pending" -> "Not Recorded"), and
#    recolor them from the "Pending" yellow to the "Not Recorded" red.
# ---------------------------------------------------------------------
$overdueRows = @(27,53,79)
foreach ($r in $overdueRows) {
    $ws.Cells.Item($r, 9).Value = "Not Recorded"
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 9)).Interior.Color = 12695295
}

# ---------------------------------------------------------------------
# 5) Widen the "Status" column now that it holds "Not Recorded"
# ---------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 13.17
